$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns L (Nb nouveaux décès à l'hôpital) and M (Nb nouveaux décès
# extra-hospitaliers) are formatted as Text (@) in this sheet. Assigning a
# plain number through .Value on a Text-formatted cell makes Excel store it
# as a text string instead of a number. To keep these as genuine numbers
# (matching the source data), temporarily switch the cell to a numeric
# format, write the value, then restore the original Text format.
function Set-NumericValue($range, $value) {
    $range.NumberFormat = "General"
    $range.Value = $value
    $range.NumberFormat = "@"
}

# --- Rows with only "Nb nouveaux cas positifs" (C) updates ---
$ws.Range("C251").Value = 850
$ws.Range("C252").Value = 682
$ws.Range("C279").Value = 164

# --- Rows 282-286: corrections to new positive cases (C), new hospital
#     deaths (L) and new extra-hospital deaths (M) ---
Set-NumericValue $ws.Range("M282") 3

$ws.Range("C283").Value = 138
Set-NumericValue $ws.Range("M283") 2

$ws.Range("C284").Value = 62
Set-NumericValue $ws.Range("M284") 2

$ws.Range("C285").Value = 63
Set-NumericValue $ws.Range("L285") 5

$ws.Range("C286").Value = 168
Set-NumericValue $ws.Range("L286") 1
Set-NumericValue $ws.Range("M286") 2

# --- Rows 287-288: newly filled-in daily figures ---
$ws.Range("C287").Value = 65
$ws.Range("E287").Value = 21
$ws.Range("F287").Value = 16
$ws.Range("G287").Value = 109
Set-NumericValue $ws.Range("L287") 1
Set-NumericValue $ws.Range("M287") 0

$ws.Range("C288").Value = 14
$ws.Range("E288").Value = 19
$ws.Range("F288").Value = 16
$ws.Range("G288").Value = 109
Set-NumericValue $ws.Range("L288") 1
Set-NumericValue $ws.Range("M288") 0
